$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue 'D2' '66.401.84'
Set-TextValue 'E2' '  +4.37%  '
Set-TextValue 'D3' '3.491.82'
Set-TextValue 'E3' '  +2.76%  '
Set-TextValue 'E4' '  -0.03%  '
Set-TextValue 'D5' '597.86'
Set-TextValue 'E5' '  +5.48%  '
Set-TextValue 'D6' '170.52'
Set-TextValue 'E6' '  +8.74%  '
Set-TextValue 'D7' '1.00'
Set-TextValue 'E7' '  -0.04%  '
Set-TextValue 'D8' '3.494.52'
Set-TextValue 'E8' '  +2.77%  '
Set-TextValue 'D9' '0.570'
Set-TextValue 'E9' '  +0.33%  '
Set-TextValue 'D10' '7.28'
Set-TextValue 'E10' '  +1.20%  '
Set-TextValue 'E11' '  +4.82%  '
Set-TextValue 'D12' '0.437'
Set-TextValue 'E12' '  +2.22%  '
Set-TextValue 'D13' '4.090.03'
Set-TextValue 'E13' '  +2.66%  '
Set-TextValue 'E14' '  +0.85%  '
Set-TextValue 'D15' '27.91'
Set-TextValue 'E15' '  +3.57%  '
Set-TextValue 'D16' '0.0000178'
Set-TextValue 'E16' '  +3.13%  '
Set-TextValue 'D17' '66.324.67'
Set-TextValue 'E17' '  +4.15%  '
Set-TextValue 'D18' '3.489.18'
Set-TextValue 'E18' '  +0.98%  '
Set-TextValue 'D19' '6.29'
Set-TextValue 'E19' '  +2.94%  '
Set-TextValue 'D20' '14.06'
Set-TextValue 'E20' '  +3.85%  '
Set-TextValue 'D21' '388.96'
Set-TextValue 'E21' '  +3.52%  '
Set-TextValue 'D22' '8.03'
Set-TextValue 'E22' '  +4.37%  '
Set-TextValue 'B23' 'Litecoin'
Set-TextValue 'C23' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D23' '72.88'
Set-TextValue 'E23' '  +2.86%  '
Set-TextValue 'B24' 'Dai'
Set-TextValue 'C24' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D24' '1.00'
Set-TextValue 'E24' '  +0.08%  '
Set-TextValue 'D25' '0.527'
Set-TextValue 'E25' '  +1.98%  '
Set-TextValue 'E26' '  +7.16%  '
Set-TextValue 'D27' '10.13'
Set-TextValue 'E27' '  +4.96%  '
Set-TextValue 'E28' '  +2.53%  '
Set-TextValue 'D29' '0.997'
Set-TextValue 'E29' '  -0.44%  '
Set-TextValue 'D30' '6.37'
Set-TextValue 'E30' '  +6.23%  '
Set-TextValue 'E31' '  +6.42%  '
Set-TextValue 'E32' '  +5.34%  '
Set-TextValue 'D33' '23.49'
Set-TextValue 'E33' '  +3.23%  '
Set-TextValue 'D34' '7.42'
Set-TextValue 'E34' '  +7.08%  '
Set-TextValue 'E35' '  +0.04%  '
Set-TextValue 'E36' '  +1.08%  '
Set-TextValue 'D37' '160.50'
Set-TextValue 'E37' '  +0.18%  '
Set-TextValue 'D38' '0.904'
Set-TextValue 'E38' '  +10.14%  '
Set-TextValue 'D39' '1.93'
Set-TextValue 'E39' '  +5.81%  '
Set-TextValue 'E40' '  +3.34%  '
Set-TextValue 'D41' '26.42'
Set-TextValue 'E41' '  +2.05%  '
Set-TextValue 'D42' '6.69'
Set-TextValue 'E42' '  +5.45%  '
Set-TextValue 'D43' '2.830.54'
Set-TextValue 'E43' '  +1.54%  '
Set-TextValue 'D44' '27.08'
Set-TextValue 'E44' '  +6.45%  '
Set-TextValue 'B45' 'Filecoin'
Set-TextValue 'C45' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D45' '4.56'
Set-TextValue 'E45' '  +3.76%  '
Set-TextValue 'B46' 'OKB'
Set-TextValue 'C46' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D46' '43.30'
Set-TextValue 'E46' '  +1.72%  '
Set-TextValue 'B47' 'VeChain'
Set-TextValue 'C47' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D47' '0.0314'
Set-TextValue 'E47' '  +3.53%  '
Set-TextValue 'B48' 'dogwifhat'
Set-TextValue 'C48' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D48' '2.52'
Set-TextValue 'E48' '  +8.09%  '
Set-TextValue 'D49' '350.00'
Set-TextValue 'E49' '  +8.72%  '
Set-TextValue 'D50' '1.10'
Set-TextValue 'E50' '  +6.71%  '
Set-TextValue 'E51' '  +9.75%  '
